$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "28.797.13"
$ws.Range("E2").Value = "  +2.88%  "
$ws.Range("D3").Value = "1.880.58"
$ws.Range("E3").Value = "  +3.14%  "
$ws.Range("E4").Value = "  +0.55%  "
Set-TextValue $ws.Range("D5") "323.58"
$ws.Range("E5").Value = "  -1.35%  "
Set-TextValue $ws.Range("D6") "1.003"
$ws.Range("E6").Value = "  +0.37%  "
Set-TextValue $ws.Range("D7") "0.4673"
$ws.Range("E7").Value = "  +1.05%  "
Set-TextValue $ws.Range("D8") "0.3932"
$ws.Range("E8").Value = "  +2.15%  "
$ws.Range("E9").Value = "  +0.87%  "
Set-TextValue $ws.Range("D10") "0.9799"
$ws.Range("E10").Value = "  +2.32%  "
Set-TextValue $ws.Range("D11") "22.31"
$ws.Range("E11").Value = "  +2.22%  "
$ws.Range("D12").Value = "1.945.12"
$ws.Range("E12").Value = "  +8.62%  "
Set-TextValue $ws.Range("D13") "5.744"
$ws.Range("E13").Value = "  +1.77%  "
Set-TextValue $ws.Range("D14") "7.007"
$ws.Range("E14").Value = "  +2.08%  "
$ws.Range("E15").Value = "  +1.64%  "
Set-TextValue $ws.Range("D16") "88.70"
$ws.Range("E16").Value = "  +2.86%  "
Set-TextValue $ws.Range("D17") "1.005"
$ws.Range("E17").Value = "  +0.42%  "
Set-TextValue $ws.Range("D18") "0.00001009"
$ws.Range("E18").Value = "  +1.65%  "
Set-TextValue $ws.Range("D19") "16.96"
$ws.Range("E19").Value = "  +2.00%  "
$ws.Range("E20").Value = "  +0.45%  "
$ws.Range("D21").Value = "28.805.01"
$ws.Range("E21").Value = "  +2.85%  "
Set-TextValue $ws.Range("D22") "5.351"
$ws.Range("E22").Value = "  +0.82%  "
Set-TextValue $ws.Range("D23") "11.09"
$ws.Range("E23").Value = "  +1.26%  "
Set-TextValue $ws.Range("D24") "2.119"
$ws.Range("E24").Value = "  +1.04%  "
$ws.Range("D25").Value = "2.123.23"
$ws.Range("E25").Value = "  +5.13%  "
Set-TextValue $ws.Range("D26") "153.61"
$ws.Range("E26").Value = "  +0.94%  "
Set-TextValue $ws.Range("D27") "19.41"
$ws.Range("E27").Value = "  +1.54%  "
Set-TextValue $ws.Range("D28") "5.760"
$ws.Range("E28").Value = "  +0.51%  "
Set-TextValue $ws.Range("D29") "2.003"
$ws.Range("E29").Value = "  +1.85%  "
Set-TextValue $ws.Range("D30") "119.95"
$ws.Range("E30").Value = "  +3.02%  "
Set-TextValue $ws.Range("D31") "0.09401"
$ws.Range("E31").Value = "  +1.97%  "
Set-TextValue $ws.Range("D32") "0.9404"
$ws.Range("E32").Value = "  +0.77%  "
Set-TextValue $ws.Range("D33") "5.316"
$ws.Range("E33").Value = "  +0.71%  "
$ws.Range("E34").Value = "  +3.22%  "
Set-TextValue $ws.Range("D35") "3.346"
$ws.Range("E35").Value = "  +0.18%  "
Set-TextValue $ws.Range("D36") "0.05916"
$ws.Range("E36").Value = "  -0.17%  "
Set-TextValue $ws.Range("D37") "0.02121"
$ws.Range("E37").Value = "  -0.75%  "
Set-TextValue $ws.Range("D38") "1.156"
$ws.Range("E38").Value = "  +1.23%  "
Set-TextValue $ws.Range("D39") "7.902"
$ws.Range("E39").Value = "  +4.22%  "
Set-TextValue $ws.Range("D40") "0.5724"
$ws.Range("E40").Value = "  +2.86%  "
Set-TextValue $ws.Range("D41") "0.1796"
$ws.Range("E41").Value = "  +1.82%  "
Set-TextValue $ws.Range("D43") "0.07301"
$ws.Range("E43").Value = "  +4.37%  "
Set-TextValue $ws.Range("D44") "11.91"
$ws.Range("E44").Value = "  +3.06%  "
Set-TextValue $ws.Range("D45") "0.5345"
$ws.Range("E45").Value = "  +2.10%  "
Set-TextValue $ws.Range("D46") "1.164"
$ws.Range("E46").Value = "  -5.25%  "
Set-TextValue $ws.Range("D47") "2.129"
$ws.Range("E47").Value = "  -3.17%  "
Set-TextValue $ws.Range("D48") "1.849"
$ws.Range("E48").Value = "  +1.54%  "
Set-TextValue $ws.Range("D49") "114.12"
$ws.Range("E49").Value = "  +2.32%  "
Set-TextValue $ws.Range("D50") "2.373"
$ws.Range("E50").Value = "  +3.02%  "
Set-TextValue $ws.Range("D51") "1.004"
$ws.Range("E51").Value = "  +0.47%  "
